# Incorporación de CHPs en el modelo y energía térmica no suministrada
#
# 1) Restructure the "Boilers" sheet: add a second product column
#    (Vitodens 020-W) before the existing Vitodens 050-W column, rename
#    "P_out" -> "P_th_nom" and "P_min" -> "P_min_porc", and drop the old
#    standalone efficiency ("n") row in favour of the new y_n / lamd_n rows.
# 2) Add a brand-new "CHPs" worksheet (after "Boilers") describing two CHP
#    unit types with nominal/electrical/thermal efficiency data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Boilers sheet: rebuild as a 3-column (A label / B / C product) table
# ---------------------------------------------------------------------
$boilers = $wb.Worksheets.Item("Boilers")

# Clear out the old 2-column table before laying out the new one.
$boilers.Cells.Clear() | Out-Null

$boilersData = @(
    @("ID",          "Vitodens 020-W",     "Vitodens 050-W"),
    @("fuel",         "Natural gas, LPG",   "Natural gas, LPG"),
    @("P_th_nom",     20,                   32),
    @("P_min_porc",   0.2,                  0.2),
    @("C_inst",       312,                  500),
    @("C_OM_kWh",     0.013,                0.013),
    @("ty",           20,                   20),
    @("y_n",          0.4576,               0.4576),
    @("lamd_n",       0.6599,               0.6599)
)

for ($i = 0; $i -lt $boilersData.Count; $i++) {
    $r = $i + 1
    $row = $boilersData[$i]
    $boilers.Cells.Item($r, 1).Value = $row[0]
    $boilers.Cells.Item($r, 2).Value = $row[1]
    $boilers.Cells.Item($r, 3).Value = $row[2]
}

$boilers.Columns.Item(2).ColumnWidth = 14.15
$boilers.Columns.Item(3).ColumnWidth = 14.3

# ---------------------------------------------------------------------
# 2. New "CHPs" sheet, placed right after "Boilers"
# ---------------------------------------------------------------------
$chps = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $boilers)
$chps.Name = "CHPs"

$chpsData = @(
    @("ID",          "Type 1 - 10kW", "Type 2 - 20kW"),
    @("C_inst",       8000,            14000),
    @("P_nom",        10,              20),
    @("P_min_porc",   0.5,             0.5),
    @("C_OM_kWh",     0.025,           0.025),
    @("ty",           15,              15),
    @("n_nom_th",     0.479,           0.486),
    @("n_nom_el",     0.405,           0.417),
    @("y_n_el",       0.2548,          0.2649),
    @("lamd_n_el",    2.2135,          2.223),
    @("y_n_th",       0.2244,          0.2355),
    @("lamd_n_th",    0.9545,          0.9655)
)

for ($i = 0; $i -lt $chpsData.Count; $i++) {
    $r = $i + 1
    $row = $chpsData[$i]
    $chps.Cells.Item($r, 1).Value = $row[0]
    $chps.Cells.Item($r, 2).Value = $row[1]
    $chps.Cells.Item($r, 3).Value = $row[2]
}

$chps.Columns.Item(2).ColumnWidth = 12.3
$chps.Columns.Item(3).ColumnWidth = 12.3

# ---------------------------------------------------------------------
# 3. Selections / active sheet bookkeeping (matches the saved view state)
# ---------------------------------------------------------------------
$boilers.Activate()
$boilers.Range("A6").Select()

$chps.Activate()
$chps.Range("D17").Select()

Write-Output "Applied Boilers restructuring and added CHPs sheet."
